# Additional companies sent for questionaire
# Remove the "Parent company" column (B) and the "Location County/City"
# column (originally E, now D after the first delete) from the known
# locomotive list, shifting the remaining columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Parent company")
$ws.Columns.Item(2).EntireColumn.Delete() | Out-Null

# Delete what is now column D ("Location County/City", originally E)
$ws.Columns.Item(4).EntireColumn.Delete() | Out-Null

# Match the author's new selection left on the data rows
$ws.Range("A2:I4").Select() | Out-Null
